$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the three new rows of progress data
$ws.Range("A10").Value = 12.9
$ws.Range("B10").Value = "黃偉倫"
$ws.Range("C10").Value = "看paper如何利用SVR技巧去把量化資料做成VIX"
$ws.Range("D10").Value = "code完成"
$ws.Range("E10").Value = "調整、篩選資料"
$ws.Range("F10").Value = "美化資料訓練模型的成果與檢查細節問題"

$ws.Range("A11").Value = 12.16
$ws.Range("B11").Value = "黃偉倫"
$ws.Range("C11").Value = "看paper如何利用SVR技巧去把量化資料做成VIX"
$ws.Range("D11").Value = "code完成"
$ws.Range("E11").Value = "模型結果符合預期"
$ws.Range("F11").Value = "完成書面報告與影片"

$ws.Range("A12").Value = 12.23
$ws.Range("B12").Value = "黃偉倫"
$ws.Range("C12").Value = "看paper如何利用SVR技巧去把量化資料做成VIX"
$ws.Range("D12").Value = "code完成"
$ws.Range("E12").Value = "全部完成"
$ws.Range("F12").Value = "全部完成"

# Update the view: move the selection to C20 (also clears the old topLeftCell scroll)
$ws.Range("C20").Select()
